$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 updates (facility_id 1002036)
$ws.Range("K3").Value = -1.819444444444444
$ws.Range("R3").Value = 1.181188907392658
$ws.Range("S3").Value = 1.232403941639294

# Row 4 updates (facility_id 1002674)
$ws.Range("K4").Value = 0.2777777777777778
$ws.Range("R4").Value = 1.194245973645681
$ws.Range("S4").Value = 1.246785162287481

# Row 8 updates (facility_id 1011252)
$ws.Range("K8").Value = 19.79629629629628
$ws.Range("R8").Value = 1.331198999020781
$ws.Range("S8").Value = 1.39868801294648

# Row 9 updates (facility_id 1013683)
$ws.Range("K9").Value = 21.28240740740739
$ws.Range("R9").Value = 1.342924567132234
$ws.Range("S9").Value = 1.411784266254412
